# Update Name of Algo
# Applies the KNN imputation result updates to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value  = 5.216
$ws.Range("A9").Value  = -21.723
$ws.Range("B12").Value = 5.57
$ws.Range("E15").Value = 16.417
$ws.Range("A18").Value = -21.997
$ws.Range("A20").Value = -20.181
$ws.Range("B26").Value = 5.705
$ws.Range("A27").Value = -21.188
$ws.Range("B27").Value = 5.695000000000001
$ws.Range("B29").Value = 5.669
$ws.Range("B37").Value = 8.376000000000001
$ws.Range("B38").Value = 5.404000000000001
$ws.Range("E38").Value = 16.557
$ws.Range("E44").Value = 16.475
$ws.Range("B51").Value = 5.603
$ws.Range("E51").Value = 16.817
$ws.Range("B55").Value = 5.705
$ws.Range("E57").Value = 16.532
$ws.Range("E63").Value = 17.673
$ws.Range("A69").Value = -21.565
$ws.Range("B69").Value = 5.669
$ws.Range("B70").Value = 5.140000000000001
$ws.Range("E70").Value = 17.523
$ws.Range("A76").Value = -20.396
$ws.Range("A82").Value = -21.993
$ws.Range("B83").Value = 5.695
$ws.Range("E99").Value = 16.583
$ws.Range("B102").Value = 7.231

$wb.Save()
